$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name, Link, Volume%) -- safe to set directly as these
# values never parse as Excel literals (numbers/dates), so .Value keeps them as text.
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.65%  '
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  -1.54%  '
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('E10').Value = '  +1.29%  '
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('E16').Value = '  -4.36%  '
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('E20').Value = '  +2.34%  '
$ws.Range('E21').Value = '  +1.11%  '
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('E30').Value = '  +1.44%  '
$ws.Range('E31').Value = '  +1.67%  '
$ws.Range('E32').Value = '  +3.80%  '
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('E35').Value = '  -0.88%  '
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('E39').Value = '  -1.48%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E48').Value = '  -2.81%  '
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('E51').Value = '  +1.48%  '

# Price column (D) -- many of these strings ARE valid Excel numeric literals
# (e.g. "1.003"), so a plain .Value assignment would silently convert them to
# numbers and reformat them (losing the exact original text). Force a Text
# number format before the write, then clear the format again afterwards so the
# cell keeps its original (default/General, unstyled) appearance -- only the
# stored value changes, matching the diff exactly.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.176.08'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.841.50'
$ws.Range('D3').ClearFormats()
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.34'
$ws.Range('D5').ClearFormats()
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').ClearFormats()
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07507'
$ws.Range('D8').ClearFormats()
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2934'
$ws.Range('D9').ClearFormats()
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '23.27'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07720'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.852.77'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.017'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6745'
$ws.Range('D14').ClearFormats()
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '83.06'
$ws.Range('D15').ClearFormats()
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.000009299'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.963'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '29.194.85'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '2.126.05'
$ws.Range('D19').ClearFormats()
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '231.82'
$ws.Range('D20').ClearFormats()
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.70'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.005'
$ws.Range('D22').ClearFormats()
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.004'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '160.51'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.546'
$ws.Range('D26').ClearFormats()
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1391'
$ws.Range('D27').ClearFormats()
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.91'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.509'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.182'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.151'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05584'
$ws.Range('D32').ClearFormats()
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.205'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7516'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.847'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.144'
$ws.Range('D36').ClearFormats()
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.667'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.778'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.225.72'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01787'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.522'
$ws.Range('D41').ClearFormats()
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.9005'
$ws.Range('D42').ClearFormats()
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.004'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.032.06'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '102.12'
$ws.Range('D45').ClearFormats()
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '66.01'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5107'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.00000000120'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.4088'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.159'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05849'
$ws.Range('D51').ClearFormats()
